# CP-159, CP-230: Using the Student to Test mapping CSV file, the StudentFactory
# creates all IRP Students to use for the simulation. They all take their
# tests synchronously.
#
# Rename the student names in column H (Student Name) from "Student<X>" to
# "IRPStudent<X>" for every data row, fix the "ARP39990002" typo in G6 to
# "AIRP39990002", and update the view's selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data rows are 3..30; column H holds the student name, which currently reads
# "StudentA".."StudentF" and must become "IRPStudentA".."IRPStudentF".
for ($row = 3; $row -le 30; $row++) {
    $cell = $ws.Cells.Item($row, 8)  # column H
    $text = $cell.Text
    if ($text -and $text.StartsWith("Student")) {
        $cell.Value = "IRP" + $text
    }
}

# Fix the typo'd AlternateSSID on row 6 (was "ARP39990002").
$ws.Range("G6").Value = "AIRP39990002"

# Update the saved view state: scroll so column C is the first visible
# column, and move the active selection to G6.
$ws.Range("G6").Select()
$excel.ActiveWindow.ScrollColumn = 3
